$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 378
$ws.Range("J2").Value = 1594
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 404
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 272
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 17
$ws.Range("S2").Value = 159
$ws.Range("T2").Value = 327
$ws.Range("U2").Value = 20
$ws.Range("V2").Value = 2480
$ws.Range("X2").Value = 2491
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 35
